$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BECbIC")

$ws.Range("B2").Value = 334498000
$ws.Range("C2").Value = 30498000
$ws.Range("D2").Value = 455922000
$ws.Range("E2").Value = 60905000
$ws.Range("F2").Value = 2524372000
$ws.Range("G2").Value = 385940000
$ws.Range("H2").Value = 755317000
$ws.Range("I2").Value = 1249502000
$ws.Range("J2").Value = 57723000
$ws.Range("K2").Value = 1620819000
$ws.Range("L2").Value = 1198491000
$ws.Range("M2").Value = 588329000
$ws.Range("N2").Value = 330884000
$ws.Range("O2").Value = 1381440000
$ws.Range("P2").Value = 1648262000
$ws.Range("Q2").Value = 748548000
$ws.Range("R2").Value = 1164066000
$ws.Range("S2").Value = 1899830000
$ws.Range("T2").Value = 1899830000
$ws.Range("U2").Value = 696553000
$ws.Range("V2").Value = 1750195000
$ws.Range("W2").Value = 14034882000
$ws.Range("X2").Value = 25994425000
$ws.Range("Y2").Value = 8380872000
$ws.Range("Z2").Value = 9204917000
$ws.Range("AA2").Value = 4292259000
$ws.Range("AB2").Value = 3419210000
$ws.Range("AC2").Value = 1605057000
$ws.Range("AD2").Value = 17367732000
$ws.Range("AE2").Value = 3789368000
$ws.Range("AF2").Value = 9936572000
$ws.Range("AG2").Value = 78479300000
$ws.Range("AH2").Value = 3888187000
$ws.Range("AI2").Value = 28121793000
$ws.Range("AJ2").Value = 2028125000
$ws.Range("AK2").Value = 616910000
